# Fruta / hortaliza, semanal
# The underlying data rows (2-19, 22-26) get shuffled: each target row receives
# the D,I,J,K,L,M,O,P values that used to belong to a different source row.
# Columns A,B,C,E,F,G,H,N,Q,R are identical for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row number -> source row number (values to copy FROM source INTO target)
$rowMap = @{
    2  = 25
    3  = 26
    4  = 19
    5  = 15
    6  = 16
    7  = 18
    8  = 6
    9  = 4
    10 = 2
    11 = 17
    12 = 22
    13 = 10
    14 = 24
    15 = 7
    16 = 11
    17 = 12
    18 = 9
    19 = 8
    22 = 23
    23 = 5
    24 = 3
    25 = 13
    26 = 14
}

# Columns (by index) whose values move together with the row permutation.
# D=4, I=9, J=10, K=11, L=12, M=13, O=15, P=16
$cols = 4,9,10,11,12,13,15,16

# First, snapshot all the original values (Value2, which gives raw numbers/strings)
# for every source row, before any writes happen.
$snapshot = @{}
foreach ($r in ($rowMap.Values | Sort-Object -Unique)) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the snapshotted values into their new target rows.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $rowVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $rowVals[$c]
    }
}
